$p = $ppt.ActivePresentation

# Add a new 5th slide using the "Title Only" layout (ppLayoutTitleOnly = 11)
$s = $p.Slides.Add(5, 11)

$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Here's a high level end-to-end architecture diagram of the Web API using the MLModelEngine with the Object Pool of PredictionEngine objects."
